$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new header columns I1 ("I0") and J1 ("IF"), copying the
# existing header cell formatting (bold, bordered, centered) from H1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for columns I (I0) and J (IF), rows 2-36.
$data = @(
  @(8, 8),
  @(7, 7),
  @(6, 6),
  @(11, 11),
  @(9, 9),
  @(8, 8),
  @(7, 7),
  @(8, 8),
  @(7, 7),
  @(7, 7),
  @(8, 8),
  @(7, 7),
  @(6, 6),
  @(9, 9),
  @(5, 6),
  @(7, 7),
  @(8, 8),
  @(8, 8),
  @(8, 8),
  @(6, 6),
  @(8, 8),
  @(9, 9),
  @(9, 9),
  @(8, 8),
  @(7, 7),
  @(8, 8),
  @(8, 8),
  @(6, 6),
  @(8, 8),
  @(8, 8),
  @(7, 7),
  @(6, 6),
  @(4, 4),
  @(7, 7),
  @(5, 5)
)

for ($i = 0; $i -lt $data.Count; $i++) {
  $row = $i + 2
  $ws.Range("I$row").Value = $data[$i][0]
  $ws.Range("J$row").Value = $data[$i][1]
}
